$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "2021" column (M) that mirrors the existing per-year columns ---

# Header cell M2 (blank, thin-bottom style) mirrors L2.
$ws.Range("L2").Copy() | Out-Null
$ws.Range("M2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Year header M3 = 2021, formatted like the other year header cells (F3's xf - no fill).
$ws.Range("F3").Copy() | Out-Null
$ws.Range("M3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("M3").Value = 2021

# Data rows: new figures for 2021 in column M.
$ws.Range("M4").Value = 7105
$ws.Range("M5").Value = 81079
$ws.Range("M6").Value = 214139

# --- Re-balance which of the two numeric-column "groups" (D:G vs H:L) carries
#     the fill-style xf so the whole D:M block again alternates consistently,
#     and drop the unused custom "0.0" number format in favour of General,
#     matching the rest of the workbook's numeric cells. ---

# Rows 4 & 5 (no outside border): D:G switch to the plain (no-fill) style,
# H:L (and the new M) switch to the fill style.
$ws.Range("B4").Copy() | Out-Null
$ws.Range("D4:G5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A4").Copy() | Out-Null
$ws.Range("H4:L5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A4").Copy() | Out-Null
$ws.Range("M4:M5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 6 (bottom border): same swap, but using the bordered donor cells.
$ws.Range("B6").Copy() | Out-Null
$ws.Range("D6:G6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A6").Copy() | Out-Null
$ws.Range("H6:L6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A6").Copy() | Out-Null
$ws.Range("M6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-apply the values (PasteSpecial of formats only shouldn't touch them, but
# make sure nothing was clobbered).
$ws.Range("M4").Value = 7105
$ws.Range("M5").Value = 81079
$ws.Range("M6").Value = 214139

# Reset the selection back to the top-left cell (matches a freshly opened view
# with nothing special selected).
$ws.Range("A1").Select() | Out-Null
